# Fruta / hortaliza, semanal
#
# Insert two new weekly price rows (Mandarina - Murcott, Terminal
# Hortofrutícola Agro Chillán, "$/caja 18 kilos") above the existing
# row 163, shifting the current rows 163:181 down to 165:183.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data (rows 163 downward) down by two rows.
$ws.Range("A163:T164").EntireRow.Insert()

# --- New row 163 ---------------------------------------------------------
$ws.Range("A163").Value = 7
$ws.Range("B163").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C163").Value = "Ñuble"
$ws.Range("D163").Value = 44783
$ws.Range("E163").Value = 16
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100102
$ws.Range("H163").Value = "Cítricos"
$ws.Range("I163").Value = 100102004
$ws.Range("J163").Value = "Mandarina"
$ws.Range("K163").Value = "Murcott"
$ws.Range("L163").Value = "Primera"
$ws.Range("M163").Value = 120
$ws.Range("N163").Value = 7500
$ws.Range("O163").Value = 8000
$ws.Range("P163").Value = 7750
$ws.Range("Q163").Value = "`$/caja 18 kilos"
$ws.Range("R163").Value = "Región de O'Higgins"
$ws.Range("S163").Value = 431
$ws.Range("T163").Value = 18

# --- New row 164 ---------------------------------------------------------
$ws.Range("A164").Value = 7
$ws.Range("B164").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C164").Value = "Ñuble"
$ws.Range("D164").Value = 44783
$ws.Range("E164").Value = 16
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100102
$ws.Range("H164").Value = "Cítricos"
$ws.Range("I164").Value = 100102004
$ws.Range("J164").Value = "Mandarina"
$ws.Range("K164").Value = "Murcott"
$ws.Range("L164").Value = "Segunda"
$ws.Range("M164").Value = 120
$ws.Range("N164").Value = 6500
$ws.Range("O164").Value = 7000
$ws.Range("P164").Value = 6750
$ws.Range("Q164").Value = "`$/caja 18 kilos"
$ws.Range("R164").Value = "Región de O'Higgins"
$ws.Range("S164").Value = 375
$ws.Range("T164").Value = 18
